# Re-applies the hourly cryptocurrency price/volume refresh described by the
# commit "Updated cryptos list on Tue Nov 19 06:42:36 UTC 2024 with GitHub Actions".
# Column D (Price) and column E (Volume(1h)) are stored as plain text in the
# sheet, and for several rows the underlying coin (columns B/C) shifted to a
# different ranking position, so whole rows of B/C/D/E are rewritten below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''91.908.00'
$ws.Range('E2').Value = '''  +0.78%  '
$ws.Range('D3').Value = '''3.129.81'
$ws.Range('E3').Value = '''  -0.45%  '
$ws.Range('D5').Value = '''241.84'
$ws.Range('E5').Value = '''  -0.25%  '
$ws.Range('D6').Value = '''618.33'
$ws.Range('E6').Value = '''  -1.49%  '
$ws.Range('E7').Value = '''  -4.85%  '
$ws.Range('D8').Value = '''0.389'
$ws.Range('E8').Value = '''  +3.94%  '
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '''  -0.05%  '
$ws.Range('D10').Value = '''3.126.44'
$ws.Range('E10').Value = '''  -0.46%  '
$ws.Range('D11').Value = '''0.754'
$ws.Range('E11').Value = '''  -1.80%  '
$ws.Range('D12').Value = '''0.206'
$ws.Range('E12').Value = '''  +0.09%  '
$ws.Range('D13').Value = '''0.0000253'
$ws.Range('E13').Value = '''  -0.14%  '
$ws.Range('D14').Value = '''35.35'
$ws.Range('E14').Value = '''  -1.60%  '
$ws.Range('E15').Value = '''  +1.45%  '
$ws.Range('D16').Value = '''91.518.48'
$ws.Range('D17').Value = '''3.709.92'
$ws.Range('E17').Value = '''  -0.34%  '
$ws.Range('D18').Value = '''3.120.39'
$ws.Range('E18').Value = '''  +0.09%  '
$ws.Range('E19').Value = '''  -0.09%  '
$ws.Range('D20').Value = '''15.02'
$ws.Range('E20').Value = '''  +1.80%  '
$ws.Range('D21').Value = '''5.93'
$ws.Range('E21').Value = '''  +0.22%  '
$ws.Range('D22').Value = '''458.43'
$ws.Range('E22').Value = '''  +1.23%  '
$ws.Range('E23').Value = '''  -5.69%  '
$ws.Range('D24').Value = '''9.28'
$ws.Range('E24').Value = '''  +1.02%  '
$ws.Range('D25').Value = '''5.95'
$ws.Range('E25').Value = '''  -0.86%  '
$ws.Range('D26').Value = '''89.60'
$ws.Range('E26').Value = '''  -4.43%  '
$ws.Range('B27').Value = '''Aptos'
$ws.Range('C27').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '''11.78'
$ws.Range('E27').Value = '''  -2.05%  '
$ws.Range('B28').Value = '''WrappedeETH'
$ws.Range('C28').Value = '''https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '''3.291.92'
$ws.Range('E28').Value = '''  +0.22%  '
$ws.Range('B29').Value = '''Hedera'
$ws.Range('C29').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '''0.145'
$ws.Range('E29').Value = '''  +17.60%  '
$ws.Range('B30').Value = '''Dai'
$ws.Range('C30').Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '''  -0.12%  '
$ws.Range('B31').Value = '''Stellar'
$ws.Range('C31').Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '''0.227'
$ws.Range('E31').Value = '''  +0.25%  '
$ws.Range('B32').Value = '''Cronos'
$ws.Range('C32').Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').Value = '''0.168'
$ws.Range('E32').Value = '''  -6.94%  '
$ws.Range('B33').Value = '''InternetComputer(DFINITY)'
$ws.Range('C33').Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''9.42'
$ws.Range('E33').Value = '''  +2.81%  '
$ws.Range('B34').Value = '''Kaspa'
$ws.Range('C34').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '''0.176'
$ws.Range('E34').Value = '''  +7.63%  '
$ws.Range('B35').Value = '''EthereumClassic'
$ws.Range('C35').Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '''26.58'
$ws.Range('E35').Value = '''  -1.63%  '
$ws.Range('B36').Value = '''RenderToken'
$ws.Range('C36').Value = '''https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').Value = '''7.50'
$ws.Range('E36').Value = '''  -3.19%  '
$ws.Range('B37').Value = '''PancakeSwap'
$ws.Range('C37').Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').Value = '''1.95'
$ws.Range('E37').Value = '''  +0.84%  '
$ws.Range('B38').Value = '''Bittensor'
$ws.Range('C38').Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '''492.66'
$ws.Range('E38').Value = '''  -2.15%  '
$ws.Range('B39').Value = '''MantraDAO'
$ws.Range('C39').Value = '''https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D39').Value = '''3.90'
$ws.Range('E39').Value = '''  -6.23%  '
$ws.Range('B40').Value = '''Fetch.AI'
$ws.Range('C40').Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '''1.32'
$ws.Range('E40').Value = '''  +0.38%  '
$ws.Range('B41').Value = '''PolygonEcosystemToken'
$ws.Range('C41').Value = '''https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = '''0.441'
$ws.Range('E41').Value = '''  +3.06%  '
$ws.Range('B42').Value = '''dogwifhat'
$ws.Range('C42').Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '''3.40'
$ws.Range('E42').Value = '''  -7.03%  '
$ws.Range('B43').Value = '''WhiteBITCoin'
$ws.Range('C43').Value = '''https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').Value = '''22.17'
$ws.Range('E43').Value = '''  +0.18%  '
$ws.Range('B44').Value = '''USDe'
$ws.Range('C44').Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '''  +0.01%  '
$ws.Range('B45').Value = '''Binance-PegBSC-USD'
$ws.Range('C45').Value = '''https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D45').Value = '''0.710'
$ws.Range('E45').Value = '''  -29.05%  '
$ws.Range('B46').Value = '''ARBITRUM'
$ws.Range('C46').Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '''0.712'
$ws.Range('E46').Value = '''  +1.21%  '
$ws.Range('B47').Value = '''Monero'
$ws.Range('C47').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '''157.39'
$ws.Range('E47').Value = '''  +0.54%  '
$ws.Range('E48').Value = '''  -0.54%  '
$ws.Range('E49').Value = '''  -0.09%  '
$ws.Range('D50').Value = '''4.49'
$ws.Range('E50').Value = '''  -2.39%  '
$ws.Range('D51').Value = '''0.0329'
$ws.Range('E51').Value = '''  +1.21%  '
